# Update marker_table genotype column (G) values: collapse homozygous
# doubled-letter genotype calls (e.g. "GG") down to a single letter (e.g. "G").
$wb = $excel.ActiveWorkbook
$markerSheet = $wb.Worksheets.Item("marker_table")

$genotypeUpdates = @{
    2  = "G"
    3  = "T"
    4  = "T"
    5  = "G"
    6  = "G"
    7  = "G"
    8  = "G"
    9  = "C"
    10 = "G"
    12 = "G"
    13 = "A"
    14 = "C"
    15 = "A"
    16 = "T"
    17 = "G"
    18 = "G"
}

foreach ($row in $genotypeUpdates.Keys) {
    $markerSheet.Range("G$row").Value = $genotypeUpdates[$row]
}

# Fill in the overall genotype result for the sample on the genotype_result sheet.
$resultSheet = $wb.Worksheets.Item("genotype_result")
$resultSheet.Range("B2").Value = "*5/*10B"
